# The two "De Havilland" rows (Comet 1, Comet 4) currently sit at the top of
# this data block (rows 52-53). The commit moves them to the bottom of the
# block (rows 73-74), shifting all the rows in between up by two, and
# re-cases the manufacturer name from "De Havilland" to "de Havilland" for
# those two relocated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the two "De Havilland" rows (52:53) to a staging area just past the
#    end of the current data block (75:76).
$ws.Range("A52:Z53").Copy()
$ws.Range("A75").PasteSpecial()

# 2. Delete the original rows 52:53. This shifts rows 54-76 up by two rows,
#    so the old rows 54-74 become rows 52-72, and the staged copy (75:76)
#    becomes rows 73-74.
$ws.Range("52:53").Delete()

# 3. Rename the manufacturer in the relocated rows from "De Havilland" to
#    "de Havilland".
$ws.Range("A73").Value = "de Havilland"
$ws.Range("A74").Value = "de Havilland"
